$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.933.27"
$ws.Range("E2").Value = "  +2.48%  "
$ws.Range("D3").Value = "2.219.62"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "263.53"
$ws.Range("E5").Value = "  +2.37%  "
$ws.Range("D6").Value = "86.82"
$ws.Range("E6").Value = "  +13.50%  "
$ws.Range("D7").Value = "0.623"
$ws.Range("E7").Value = "  +1.36%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +2.36%  "
$ws.Range("D10").Value = "46.51"
$ws.Range("E10").Value = "  +10.75%  "
$ws.Range("D11").Value = "0.0922"
$ws.Range("E11").Value = "  +2.07%  "
$ws.Range("D12").Value = "7.61"
$ws.Range("E12").Value = "  +9.31%  "
$ws.Range("E13").Value = "  +2.70%  "
$ws.Range("D14").Value = "2.551.81"
$ws.Range("E14").Value = "  +0.27%  "
$ws.Range("D15").Value = "14.68"
$ws.Range("E15").Value = "  +1.31%  "
$ws.Range("D16").Value = "2.213.94"
$ws.Range("E16").Value = "  -0.21%  "
$ws.Range("D17").Value = "0.783"
$ws.Range("E17").Value = "  +0.38%  "
$ws.Range("D18").Value = "43.921.62"
$ws.Range("E18").Value = "  +2.49%  "
$ws.Range("E19").Value = "  +1.18%  "
$ws.Range("D20").Value = "6.00"
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("D21").Value = "70.11"
$ws.Range("E21").Value = "  -1.69%  "
$ws.Range("D22").Value = "2.39"
$ws.Range("E22").Value = "  +8.54%  "
$ws.Range("D23").Value = "232.27"
$ws.Range("E23").Value = "  +0.78%  "
$ws.Range("D24").Value = "9.04"
$ws.Range("E24").Value = "  -2.80%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").Value = "10.76"
$ws.Range("E26").Value = "  -0.21%  "
$ws.Range("D27").Value = "3.52"
$ws.Range("E27").Value = "  +5.28%  "
$ws.Range("D28").Value = "39.91"
$ws.Range("E28").Value = "  -5.33%  "
$ws.Range("D29").Value = "2.27"
$ws.Range("E29").Value = "  +3.27%  "
$ws.Range("D30").Value = "2.26"
$ws.Range("E30").Value = "  +1.90%  "
$ws.Range("D31").Value = "174.99"
$ws.Range("E31").Value = "  +0.49%  "
$ws.Range("E32").Value = "  +1.34%  "
$ws.Range("D33").Value = "0.0886"
$ws.Range("E33").Value = "  +1.34%  "
$ws.Range("E34").Value = "  +4.06%  "
$ws.Range("E35").Value = "  +1.53%  "
$ws.Range("D36").Value = "0.111"
$ws.Range("E36").Value = "  +4.31%  "
$ws.Range("D37").Value = "0.0362"
$ws.Range("E37").Value = "  +0.44%  "
$ws.Range("E38").Value = "  +3.33%  "
$ws.Range("D39").Value = "3.26"
$ws.Range("E39").Value = "  +14.97%  "
$ws.Range("D40").Value = "12.46"
$ws.Range("E40").Value = "  -3.46%  "
$ws.Range("D41").Value = "65.13"
$ws.Range("E41").Value = "  +8.50%  "
$ws.Range("E42").Value = "  -0.64%  "
$ws.Range("D43").Value = "5.56"
$ws.Range("E43").Value = "  +4.70%  "
$ws.Range("D44").Value = "0.204"
$ws.Range("E44").Value = "  +2.26%  "
$ws.Range("D45").Value = "101.41"
$ws.Range("E45").Value = "  -1.29%  "
$ws.Range("E46").Value = "  +0.55%  "
$ws.Range("D47").Value = "0.0985"
$ws.Range("E47").Value = "  +0.78%  "
$ws.Range("E48").Value = "  +1.51%  "
$ws.Range("E49").Value = "  +4.22%  "
$ws.Range("D50").Value = "0.448"
$ws.Range("E50").Value = "  -2.38%  "
$ws.Range("D51").Value = "1.54"
$ws.Range("E51").Value = "  +7.57%  "